$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1433
$ws.Range("E2").Value = 808
$ws.Range("F2").Value = 808
$ws.Range("G2").Value = 802
$ws.Range("H2").Value = 599
$ws.Range("I2").Value = 598
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 5175
$ws.Range("L2").Value = 868
$ws.Range("M2").Value = 4307
$ws.Range("N2").Value = 4280
$ws.Range("O2").Value = 28
$ws.Range("P2").Value = 2525
$ws.Range("Q2").Value = 571
$ws.Range("R2").Value = 382
$ws.Range("S2").Value = -1040
$ws.Range("T2").Value = 4
$ws.Range("V2").Value = 323
$ws.Range("W2").Value = 56.39
$ws.Range("X2").Value = 41.82
$ws.Range("Y2").Value = 14.8
$ws.Range("Z2").Value = 11.16
$ws.Range("AA2").Value = 20.14
$ws.Range("AB2").Value = 75.36
$ws.Range("AC2").Value = 237
$ws.Range("AD2").Value = 15.73
$ws.Range("AE2").Value = 1745
$ws.Range("AF2").Value = 2.13
$ws.Range("AG2").Value = 40
$ws.Range("AH2").Value = 1.07
$ws.Range("AI2").Value = 16.4
$ws.Range("AJ2").Value = 252489230
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 1385
$ws.Range("E3").Value = 890
$ws.Range("F3").Value = 890
$ws.Range("G3").Value = 896
$ws.Range("H3").Value = 682
$ws.Range("I3").Value = 682
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6534
$ws.Range("L3").Value = 1672
$ws.Range("M3").Value = 4862
$ws.Range("N3").Value = 4862
$ws.Range("P3").Value = 2525
$ws.Range("Q3").Value = 534
$ws.Range("R3").Value = -1156
$ws.Range("S3").Value = 531
$ws.Range("T3").Value = 6
$ws.Range("V3").Value = 955
$ws.Range("W3").Value = 64.27
$ws.Range("X3").Value = 49.24
$ws.Range("Y3").Value = 14.92
$ws.Range("Z3").Value = 11.64
$ws.Range("AA3").Value = 34.4
$ws.Range("AB3").Value = 97.33
$ws.Range("AC3").Value = 270
$ws.Range("AD3").Value = 10.48
$ws.Range("AE3").Value = 1982
$ws.Range("AF3").Value = 1.43
$ws.Range("AG3").Value = 60
$ws.Range("AH3").Value = 2.12
$ws.Range("AI3").Value = 21.59
$ws.Range("AJ3").Value = 252489230
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 1780
$ws.Range("E4").Value = 1140
$ws.Range("F4").Value = 1140
$ws.Range("G4").Value = 1591
$ws.Range("H4").Value = 1214
$ws.Range("I4").Value = 1213
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 9544
$ws.Range("L4").Value = 3641
$ws.Range("M4").Value = 5903
$ws.Range("N4").Value = 5872
$ws.Range("O4").Value = 31
$ws.Range("P4").Value = 2525
$ws.Range("Q4").Value = -1249
$ws.Range("R4").Value = -497
$ws.Range("S4").Value = 1727
$ws.Range("T4").Value = 8
$ws.Range("V4").Value = 2797
$ws.Range("W4").Value = 64.04000000000001
$ws.Range("X4").Value = 68.23
$ws.Range("Y4").Value = 22.61
$ws.Range("Z4").Value = 15.11
$ws.Range("AA4").Value = 61.68
$ws.Range("AB4").Value = 138.56
$ws.Range("AC4").Value = 481
$ws.Range("AD4").Value = 5.94
$ws.Range("AE4").Value = 2394
$ws.Range("AF4").Value = 1.19
$ws.Range("AG4").Value = 75
$ws.Range("AH4").Value = 2.63
$ws.Range("AI4").Value = 15.16
$ws.Range("AJ4").Value = 252489230
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 2500
$ws.Range("E5").Value = 1711
$ws.Range("F5").Value = 1711
$ws.Range("G5").Value = 2210
$ws.Range("H5").Value = 1678
$ws.Range("I5").Value = 1672
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 12107
$ws.Range("L5").Value = 4915
$ws.Range("M5").Value = 7192
$ws.Range("N5").Value = 7155
$ws.Range("O5").Value = 37
$ws.Range("P5").Value = 2525
$ws.Range("Q5").Value = -966
$ws.Range("R5").Value = 371
$ws.Range("S5").Value = 794
$ws.Range("T5").Value = 4
$ws.Range("V5").Value = 3794
$ws.Range("W5").Value = 68.45999999999999
$ws.Range("X5").Value = 67.13
$ws.Range("Y5").Value = 25.66
$ws.Range("Z5").Value = 15.5
$ws.Range("AA5").Value = 68.33
$ws.Range("AB5").Value = 189.63
$ws.Range("AC5").Value = 662
$ws.Range("AD5").Value = 4.73
$ws.Range("AE5").Value = 2917
$ws.Range("AF5").Value = 1.07
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 3.19
$ws.Range("AI5").Value = 14.67
$ws.Range("AJ5").Value = 252489230
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 2693
$ws.Range("E6").Value = 1846
$ws.Range("F6").Value = 1846
$ws.Range("G6").Value = 2229
$ws.Range("H6").Value = 1670
$ws.Range("I6").Value = 1671
$ws.Range("K6").Value = 13759
$ws.Range("L6").Value = 5513
$ws.Range("M6").Value = 8246
$ws.Range("N6").Value = 8210
$ws.Range("P6").Value = 2525
$ws.Range("Q6").Value = 661
$ws.Range("R6").Value = 308
$ws.Range("S6").Value = 334
$ws.Range("T6").Value = 5
$ws.Range("V6").Value = 4494
$ws.Range("W6").Value = 68.53
$ws.Range("X6").Value = 62
$ws.Range("Y6").Value = 21.75
$ws.Range("Z6").Value = 12.91
$ws.Range("AA6").Value = 66.86
$ws.Range("AB6").Value = 233.59
$ws.Range("AC6").Value = 662
$ws.Range("AD6").Value = 4.08
$ws.Range("AE6").Value = 3379
$ws.Range("AF6").Value = 0.8
$ws.Range("AG6").Value = 110
$ws.Range("AH6").Value = 4.07
$ws.Range("AI6").Value = 15.99
$ws.Range("AJ6").Value = 252489230
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 2652
$ws.Range("E7").Value = 1509
$ws.Range("G7").Value = 1737
$ws.Range("H7").Value = 1258
$ws.Range("I7").Value = 1256
$ws.Range("K7").Value = 15162
$ws.Range("L7").Value = 5650
$ws.Range("M7").Value = 9502
$ws.Range("N7").Value = 9465
$ws.Range("P7").Value = 2522
$ws.Range("W7").Value = 56.91
$ws.Range("X7").Value = 47.46
$ws.Range("Y7").Value = 14.22
$ws.Range("Z7").Value = 8.699999999999999
$ws.Range("AA7").Value = 59.46
$ws.Range("AC7").Value = 498
$ws.Range("AD7").Value = 4.03
$ws.Range("AE7").Value = 3965
$ws.Range("AF7").Value = 0.51
$ws.Range("AG7").Value = 100
$ws.Range("AH7").Value = 4.99
$ws.Range("AI7").Value = 20.09
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 2464
$ws.Range("E8").Value = 1346
$ws.Range("G8").Value = 1560
$ws.Range("H8").Value = 1140
$ws.Range("I8").Value = 1142
$ws.Range("K8").Value = 15924
$ws.Range("L8").Value = 5485
$ws.Range("M8").Value = 10398
$ws.Range("N8").Value = 10361
$ws.Range("P8").Value = 2522
$ws.Range("W8").Value = 54.66
$ws.Range("X8").Value = 46.3
$ws.Range("Y8").Value = 11.52
$ws.Range("Z8").Value = 7.34
$ws.Range("AA8").Value = 52.75
$ws.Range("AC8").Value = 452
$ws.Range("AD8").Value = 4.43
$ws.Range("AE8").Value = 4341
$ws.Range("AF8").Value = 0.46
$ws.Range("AG8").Value = 90
$ws.Range("AH8").Value = 4.49
$ws.Range("AI8").Value = 19.91
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 2310
$ws.Range("E9").Value = 1263
$ws.Range("G9").Value = 1480
$ws.Range("H9").Value = 1085
$ws.Range("I9").Value = 1083
$ws.Range("K9").Value = 16752
$ws.Range("L9").Value = 5427
$ws.Range("M9").Value = 11260
$ws.Range("N9").Value = 11222
$ws.Range("P9").Value = 2522
$ws.Range("W9").Value = 54.69
$ws.Range("X9").Value = 46.98
$ws.Range("Y9").Value = 10.04
$ws.Range("Z9").Value = 6.64
$ws.Range("AA9").Value = 48.2
$ws.Range("AC9").Value = 429
$ws.Range("AD9").Value = 4.67
$ws.Range("AE9").Value = 4702
$ws.Range("AF9").Value = 0.43
$ws.Range("AG9").Value = 90
$ws.Range("AH9").Value = 4.49
$ws.Range("AI9").Value = 20.98
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
